# Update NTT Data's address block on the title page.
$d = $word.ActiveDocument

# 1) Company name -> add legal suffix "S.A."
$d.Content.Find.Execute("NTT Data Romania", $true, $false, $false, $false, $false, $true, 1, $false, "NTT Data Romania S.A.", 2) | Out-Null

# 2) Street line -> reorder to "<number>, <street> Street,"
$d.Content.Find.Execute("Street Constanta 19-21 ", $true, $false, $false, $false, $false, $true, 1, $false, "19-21, Constanta Street,", 2) | Out-Null

# 3) City line -> postal code now leads the city name
$d.Content.Find.Execute("Cluj Napoca City, 400158", $true, $false, $false, $false, $false, $true, 1, $false, "400158 Cluj Napoca", 2) | Out-Null
